# Aula T2-S01 e T2-S02, updates de código
# Nudge the main content textbox(es) on each slide up slightly
# (a:off y changes while a:ext stays the same).

$p = $ppt.ActivePresentation

# EMU -> points conversion (1 pt = 12700 EMU).
# Shape.Top/.Left are single-precision (Single) in the PowerPoint object
# model, and the host truncates (floors) pt*12700 back to EMU on save, so
# add a tiny (0.5 EMU) bias before narrowing to f32 to land on the exact
# target EMU value instead of one EMU short.
function EMUToPt($emu) { return ($emu + 0.5) / 12700.0 }

# (slideIndex, shapeIndex, newTopEMU)
$targets = @(
    @(2,  2, 841276),
    @(3,  2, 841276),
    @(4,  2, 841276),
    @(5,  2, 841276),
    @(6,  2, 841276),
    @(7,  2, 841276),
    @(7,  4, 1561356),
    @(8,  2, 841276),
    @(8,  4, 2353444),
    @(9,  2, 841276),
    @(10, 2, 841276),
    @(11, 2, 841276),
    @(12, 2, 841276),
    @(13, 2, 841276),
    @(14, 2, 841276),
    @(15, 2, 841276),
    @(16, 2, 841276),
    @(17, 2, 841276),
    @(18, 2, 841276)
)

foreach ($t in $targets) {
    $slideIdx = $t[0]
    $shapeIdx = $t[1]
    $newTopEmu = $t[2]

    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item($shapeIdx)
    $shape.Top = EMUToPt $newTopEmu
}
